# Add a new "GAAP GL Account Class Profile" data row to the TestData sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("TestData")

$ws.Range("A5").Value = 4

# Write the all-caps string first, then the title-case one, so the new
# shared-string table entries land in the same order as the authored file.
$ws.Range("C5").Value = " GAAP GL ACCOUNT CLASS PROFILE "
$ws.Range("B5").Value = "Gaap Gl Account Class Profile"
